# Weekly update: shift the existing Pepino dulce rows down by one week
# and insert the newly reported week at the top of the data block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the new "latest" record that used to occupy row 11 (a duplicate
#    of the former row 10 data, now one week further back) before touching
#    the existing rows, so we don't lose the values we need to copy down.
$ws.Range("A11").Value = 7
$ws.Range("B11").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C11").Value = "Ñuble"
$ws.Range("D11").Value = 44628
$ws.Range("D11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 100112043
$ws.Range("G11").Value = "Pepino dulce"
$ws.Range("H11").Value = "Cultivar IV Región"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 60
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 16000
$ws.Range("M11").Value = 15500
$ws.Range("N11").Value = "$/bandeja 18 kilos"
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 861
$ws.Range("Q11").Value = 18
$ws.Range("R11").Value = "Hortaliza"

# 2) Shift the dates/volumes of the existing weekly records down one week.
$ws.Range("D10").Value = 44637
$ws.Range("J10").Value = 100

$ws.Range("D9").Value = 44642

$ws.Range("D8").Value = 44651
$ws.Range("J8").Value = 60
